# repull data, push all data, mean calculation
# Update the dSF (column F) values for the rows whose data was repulled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = 0
    8  = 2
    10 = 0
    16 = 0
    21 = -4
    26 = -7
    28 = -7
    29 = -7
    31 = -5
    32 = -2
    33 = 5
    35 = -1
    36 = -2
    43 = 0
    47 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
